$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: remove the standalone "Spatial Transformation. A vital module..."
# paragraph together with the blank paragraph that precedes it (right after
# the "MLP Based Spatial Transformation Pyramid" section's body text).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Spatial Transformation. A vital module of 3D lane detection") | Out-Null
$spatialPara = $rng.Paragraphs(1)
$blankBefore = $spatialPara.Previous()

# Delete the body paragraph first, then the blank paragraph above it, so the
# indices of earlier content are unaffected while we do it.
$spatialPara.Range.Delete()
$blankBefore.Range.Delete()

# ---------------------------------------------------------------------------
# Change 2: the "Inference" heading loses its bold/large-size formatting
# (becomes a plain paragraph) and the big bold heading formatting, together
# with the <w:lastRenderedPageBreak/>, moves onto the "Experiments" heading
# that follows it.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Inference") | Out-Null
$inferencePara = $rng2.Paragraphs(1)
$inferencePara.Range.InsertXML("<w:p $wNs><w:r><w:t>Inference</w:t></w:r></w:p>")

$rng3 = $d.Content
$rng3.Find.Execute("Experiments") | Out-Null
$experimentsPara = $rng3.Paragraphs(1)
$experimentsPara.Range.InsertXML("<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Experiments</w:t></w:r></w:p>")

Write-Output "done"
